$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = -7.624699999999991

$ws.Range("A4").Value = -21.10420000000001
$ws.Range("C4").Value = -11.11439999999999
$ws.Range("D4").Value = -6.590299999999998

$ws.Range("C5").Value = -14.82880000000001

$ws.Range("A6").Value = -21.53370000000002

$ws.Range("A7").Value = -21.76940000000001

$ws.Range("C8").Value = -12.2887

$ws.Range("D9").Value = -7.6689

$ws.Range("D11").Value = -8.463399999999995

$ws.Range("D14").Value = -6.1902

$ws.Range("A16").Value = -20.20819999999998
$ws.Range("C16").Value = -12.00620000000001

$ws.Range("D18").Value = -8.512799999999997

$ws.Range("A20").Value = -22.79600000000002

$ws.Range("C22").Value = -11.11159999999999

$ws.Range("D25").Value = -8.06279999999999
